$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change every C cell currently valued 500 down to 50 (rows 1-16)
for ($r = 1; $r -le 16; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value() -eq 500) {
        $cell.Value = 50
    }
}

# Move the active selection to C16, matching the final cursor position after edits
$ws.Range("C16").Select()
